# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Update the "Periodo Mora" (column E) entries in the worker arrears table:
# the oldest period (2412) drops to the bottom of the recent-period block and
# the five most-recent periods (2501-2505) shift up, now listed most-recent-first.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2503"
$ws.Range("E19").Value = "2502"
$ws.Range("E20").Value = "2501"
$ws.Range("E21").Value = "2412"
$ws.Range("E22").Value = "2506"
